# Natmi following Dr Hou advice
# Update the LR-pair results sheet: recompute row 2 with the corrected
# values, and expand the result table to cover the full set of
# Sending/Target cluster combinations (FAPs/sCs x FAPs/sCs) for the
# Wnt9a -> Fzd10 ligand-receptor pair (rows 2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("FAPs", "Wnt9a", "Fzd10", "FAPs", 3, 1, 3.232908333333333, 9.698725, 0.5584514397475191, 0.5584514397475192, 3, 1, 0.06694666666666667, 0.20084, 0.5868893752684747, 0.5868893752684747, 0.2164324365555556, 1.947891929, 0.3277492165912018, 0.3277492165912018)
$row3 = @("FAPs", "Wnt9a", "Fzd10", "sCs", 3, 1, 3.232908333333333, 9.698725, 0.5584514397475191, 0.5584514397475192, 2, 0.6666666666666666, 0.04712366666666667, 0.141371, 0.4131106247315253, 0.4131106247315252, 0.1523464946638889, 1.371118451975, 0.2307022231563174, 0.2307022231563174)
$row4 = @("sCs", "Wnt9a", "Fzd10", "FAPs", 3, 1, 2.556150666666667, 7.668452, 0.4415485602524809, 0.441548560252481, 3, 1, 0.06694666666666667, 0.20084, 0.5868893752684747, 0.5868893752684747, 0.1711257666311111, 1.54013189968, 0.259140158677273, 0.259140158677273)
$row5 = @("sCs", "Wnt9a", "Fzd10", "sCs", 3, 1, 2.556150666666667, 7.668452, 0.4415485602524809, 0.441548560252481, 2, 0.6666666666666666, 0.04712366666666667, 0.141371, 0.4131106247315253, 0.4131106247315252, 0.1204551919657778, 1.084096727692, 0.1824084015752079, 0.1824084015752079)

$rowNums = @(2, 3, 4, 5)
$rowData = @($row2, $row3, $row4, $row5)

for ($ridx = 0; $ridx -lt $rowNums.Length; $ridx++) {
    $r = $rowNums[$ridx]
    $vals = $rowData[$ridx]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
